# Agenda de projeto.xlsx - update "Ide/Criação/Teste" status cells
# from "?" to "!" for the rows that are now complete, and move the
# active cell selection to F18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("F13").Value = "!"
$ws.Range("F16").Value = "!"
$ws.Range("F17").Value = "!"

$ws.Range("F18").Select()
